$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# Add the new value in C10 (fills in the previously-empty "Accuracy" cell for SDU3.3)
$ws.Range("C10").Value = 0.7890625

# Give C1 (header "Accuracy") a solid fill (new style) to highlight it
$ws.Range("C1").Interior.ThemeColor = 0
$ws.Range("C1").Interior.PatternColorIndex = 64
$ws.Range("C1").Interior.Pattern = -4124

# Move the active selection to C11
$ws.Range("C11").Select()

$wb.Save()
